$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common repeated values
$dipc = "Donostia International Physics Center, BERC Basque Excellence Research Centre, Manuel de Lardizabal 4"
$sansebastian = "San Sebasti\'an / Donostia, E-20018, Spain"

# 1. Move A.I. Aranburu's Institution1 from CFM-area entry to DIPC
$ws.Range("E6").Value = $dipc

# 2. Insert new row for Castillo, A. (alphabetically before Cebrian, currently row 22)
$ws.Rows.Item(22).Insert()
$ws.Range("A22").Value = "Castillo"
$ws.Range("B22").Value = "A.       "
$ws.Range("E22").Value = $dipc
$ws.Range("F22").Value = $sansebastian

# 3. Insert new row for Echevarria, C. (alphabetically between Dickel and Elorza, now row 32)
$ws.Rows.Item(32).Insert()
$ws.Range("A32").Value = "Echevarria"
$ws.Range("B32").Value = "C."
$ws.Range("E32").Value = $dipc
$ws.Range("F32").Value = $sansebastian

# 4. Insert new row for Seeman, M. (alphabetically between Santos and Shomroni, now row 100)
$ws.Rows.Item(100).Insert()
$ws.Range("A100").Value = "Seeman"
$ws.Range("B100").Value = "M."
$ws.Range("E100").Value = $dipc
$ws.Range("F100").Value = $sansebastian

# 5. Append new row for Yubero, A. at the very end (new last row, 114)
$ws.Range("A114").Value = "Yubero"
$ws.Range("B114").Value = "A."
$ws.Range("E114").Value = $dipc
$ws.Range("F114").Value = $sansebastian
